$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy style (s="1") from an existing column-A data cell down to the new rows 49:74
$ws.Range("A2").Copy($ws.Range("A49:A74"))

# Update existing rows 2-48 (columns B, C, D) with new values
$ws.Range("B2").Value = 1.452655467151331
$ws.Range("C2").Value = 9.002301462397778
$ws.Range("D2").Value = 0.14453685936482
$ws.Range("B3").Value = 1.530240236146303
$ws.Range("C3").Value = 5.662932838612557
$ws.Range("D3").Value = 0.1646698078461855
$ws.Range("B4").Value = 2.09693369154597
$ws.Range("C4").Value = 7.417766695608578
$ws.Range("D4").Value = 0.3603193133634053
$ws.Range("B5").Value = 2.375704881924857
$ws.Range("C5").Value = 2.651358042865527
$ws.Range("D5").Value = 0.1149751518475411
$ws.Range("B6").Value = 3.221475119261219
$ws.Range("C6").Value = 2.731300435084626
$ws.Range("D6").Value = 0.2618408025998836
$ws.Range("B7").Value = 3.735553683524884
$ws.Range("C7").Value = 4.575908499870541
$ws.Range("D7").Value = 0.1912190204376858
$ws.Range("B8").Value = 4.65617803523349
$ws.Range("C8").Value = 4.241836350494927
$ws.Range("D8").Value = 0.2589437323175452
$ws.Range("B9").Value = 6.717545901270031
$ws.Range("C9").Value = 7.285945697967334
$ws.Range("D9").Value = 0.3484315621182372
$ws.Range("B10").Value = 8.44668177676159
$ws.Range("C10").Value = 3.71962737981124
$ws.Range("D10").Value = 0.2013327805231728
$ws.Range("B11").Value = 8.540373289135248
$ws.Range("C11").Value = 6.256924445512367
$ws.Range("D11").Value = 0.3990945310545316
$ws.Range("B12").Value = 11.21653334512462
$ws.Range("C12").Value = 3.829948852636527
$ws.Range("D12").Value = 0.2301572287837012
$ws.Range("B13").Value = 15.83883632512221
$ws.Range("C13").Value = 4.319528713652799
$ws.Range("D13").Value = 0.09623868781778294
$ws.Range("B14").Value = 16.9129082917815
$ws.Range("C14").Value = 5.773541513956777
$ws.Range("D14").Value = 0.2652380919228087
$ws.Range("B15").Value = 17.38856521216392
$ws.Range("C15").Value = 7.698256076622296
$ws.Range("D15").Value = 0.4613594793625123
$ws.Range("B16").Value = 19.91536174196453
$ws.Range("C16").Value = 6.692394043266505
$ws.Range("D16").Value = 0.386762361739541
$ws.Range("B17").Value = 22.85897271718316
$ws.Range("C17").Value = 4.620255382778632
$ws.Range("D17").Value = 0.08152527407139126
$ws.Range("B18").Value = 26.67320520941897
$ws.Range("C18").Value = 3.176985746275943
$ws.Range("D18").Value = 0.2015213139249075
$ws.Range("B19").Value = 27.66354702550375
$ws.Range("C19").Value = 14.08303608642209
$ws.Range("D19").Value = 0.2149323433376239
$ws.Range("B20").Value = 28.70247727028161
$ws.Range("C20").Value = 5.313851284673716
$ws.Range("D20").Value = 0.3050444352187603
$ws.Range("B21").Value = 30.05858695471986
$ws.Range("C21").Value = 2.968920516797199
$ws.Range("D21").Value = 0.279827001873024
$ws.Range("B22").Value = 30.56893773599286
$ws.Range("C22").Value = 6.766039227958732
$ws.Range("D22").Value = 0.3366509627694773
$ws.Range("B23").Value = 31.09277546468668
$ws.Range("C23").Value = 3.656732824040027
$ws.Range("D23").Value = 0.5902953229219223
$ws.Range("B24").Value = 31.1267731636808
$ws.Range("C24").Value = 7.347709472797354
$ws.Range("D24").Value = 0.4944925523735685
$ws.Range("B25").Value = 32.27245640385294
$ws.Range("C25").Value = 5.040453215570746
$ws.Range("D25").Value = 0.3992900241488959
$ws.Range("B26").Value = 36.4937882261236
$ws.Range("C26").Value = 5.256632336748592
$ws.Range("D26").Value = 0.4464011464351865
$ws.Range("B27").Value = 37.7606065617028
$ws.Range("C27").Value = 6.302875366151429
$ws.Range("D27").Value = 0.3024330415014818
$ws.Range("B28").Value = 38.90022729922902
$ws.Range("C28").Value = 8.672912666682175
$ws.Range("D28").Value = 0.294558980864442
$ws.Range("B29").Value = 40.6600981704441
$ws.Range("C29").Value = 7.404715578667911
$ws.Range("D29").Value = 0.2997867074351855
$ws.Range("B30").Value = 42.4426235981546
$ws.Range("C30").Value = 7.535619559383209
$ws.Range("D30").Value = 0.3233352275849599
$ws.Range("B31").Value = 42.99122157416637
$ws.Range("C31").Value = 4.342673457589579
$ws.Range("D31").Value = 0.1615180413393928
$ws.Range("B32").Value = 44.00313972899117
$ws.Range("C32").Value = 5.037195133038523
$ws.Range("D32").Value = 0.1822829325935605
$ws.Range("B33").Value = 44.65068332439179
$ws.Range("C33").Value = 5.704414598422658
$ws.Range("D33").Value = 0.391175677757596
$ws.Range("B34").Value = 45.90126257486216
$ws.Range("C34").Value = 4.868958522207582
$ws.Range("D34").Value = 0.4090867395704323
$ws.Range("B35").Value = 48.51677440070429
$ws.Range("C35").Value = 7.726404403163232
$ws.Range("D35").Value = 0.293834379235222
$ws.Range("B36").Value = 48.59996311526476
$ws.Range("C36").Value = 7.400503994431939
$ws.Range("D36").Value = 0.2020852561067015
$ws.Range("B37").Value = 49.24374916834713
$ws.Range("C37").Value = 5.189907177153068
$ws.Range("D37").Value = 0.3164227573648053
$ws.Range("B38").Value = 50.13056977584841
$ws.Range("C38").Value = 8.211183853928175
$ws.Range("D38").Value = 0.1327181373012397
$ws.Range("B39").Value = 54.54366474609868
$ws.Range("C39").Value = 5.85797148604416
$ws.Range("D39").Value = 0.4296507218438935
$ws.Range("B40").Value = 56.40793420904466
$ws.Range("C40").Value = 3.008149243530818
$ws.Range("D40").Value = 0.5213511509104333
$ws.Range("B41").Value = 56.42351474086097
$ws.Range("C41").Value = 9.056637527335804
$ws.Range("D41").Value = 0.5446762026520946
$ws.Range("B42").Value = 57.01651287886276
$ws.Range("C42").Value = 7.689946999033065
$ws.Range("D42").Value = 0.2339819493093075
$ws.Range("B43").Value = 57.78606009423005
$ws.Range("C43").Value = 6.466947197068706
$ws.Range("D43").Value = 0.3722480964071053
$ws.Range("B44").Value = 58.95415454993866
$ws.Range("C44").Value = 5.784845704224262
$ws.Range("D44").Value = 0.3215796441759429
$ws.Range("B45").Value = 59.30088846109105
$ws.Range("C45").Value = 2.590796877754809
$ws.Range("D45").Value = 0.1485218900214344
$ws.Range("B46").Value = 60.41409635776665
$ws.Range("C46").Value = 5.397674160995381
$ws.Range("D46").Value = 0.4838627063395299
$ws.Range("B47").Value = 60.76777699211669
$ws.Range("C47").Value = 10.48609903963651
$ws.Range("D47").Value = 0.4515972110822398
$ws.Range("B48").Value = 61.74569767529582
$ws.Range("C48").Value = 3.993341543643208
$ws.Range("D48").Value = 0.4301064775574309

# Add new rows 49-74 (columns A, B, C, D)
$ws.Range("A49").Value = 47
$ws.Range("B49").Value = 62.00473611451248
$ws.Range("C49").Value = 9.392235147462969
$ws.Range("D49").Value = 0.443531329953228
$ws.Range("A50").Value = 48
$ws.Range("B50").Value = 65.54272817028023
$ws.Range("C50").Value = 13.60894398092449
$ws.Range("D50").Value = 0.4312612642095092
$ws.Range("A51").Value = 49
$ws.Range("B51").Value = 67.31279205083058
$ws.Range("C51").Value = 10.44467555843544
$ws.Range("D51").Value = 0.3071860441617649
$ws.Range("A52").Value = 50
$ws.Range("B52").Value = 68.14935492224568
$ws.Range("C52").Value = 8.922543025035742
$ws.Range("D52").Value = 0.139519665880157
$ws.Range("A53").Value = 51
$ws.Range("B53").Value = 68.51068354760311
$ws.Range("C53").Value = 9.963267606988172
$ws.Range("D53").Value = 0.492949237450105
$ws.Range("A54").Value = 52
$ws.Range("B54").Value = 68.85913771095731
$ws.Range("C54").Value = 10.44707811292284
$ws.Range("D54").Value = 0.3824173233751409
$ws.Range("A55").Value = 53
$ws.Range("B55").Value = 69.36338781462435
$ws.Range("C55").Value = 10.42435600266985
$ws.Range("D55").Value = 0.267847760911854
$ws.Range("A56").Value = 54
$ws.Range("B56").Value = 71.07306121371091
$ws.Range("C56").Value = 4.926013043616745
$ws.Range("D56").Value = 0.2698917643130602
$ws.Range("A57").Value = 55
$ws.Range("B57").Value = 71.82827777409398
$ws.Range("C57").Value = 4.446854631473042
$ws.Range("D57").Value = 0.1412125341062782
$ws.Range("A58").Value = 56
$ws.Range("B58").Value = 73.78774324403682
$ws.Range("C58").Value = 7.469340231789071
$ws.Range("D58").Value = 0.4850939065560733
$ws.Range("A59").Value = 57
$ws.Range("B59").Value = 76.33432613174223
$ws.Range("C59").Value = 5.549196705854865
$ws.Range("D59").Value = 0.271086965928218
$ws.Range("A60").Value = 58
$ws.Range("B60").Value = 77.01376224848245
$ws.Range("C60").Value = 3.68262041278794
$ws.Range("D60").Value = 0.1869398790624003
$ws.Range("A61").Value = 59
$ws.Range("B61").Value = 79.70323568647551
$ws.Range("C61").Value = 3.193088282135918
$ws.Range("D61").Value = 0.2414863334179861
$ws.Range("A62").Value = 60
$ws.Range("B62").Value = 80.52708713568532
$ws.Range("C62").Value = 5.916230680364833
$ws.Range("D62").Value = 0.1094226847938124
$ws.Range("A63").Value = 61
$ws.Range("B63").Value = 84.42186598340606
$ws.Range("C63").Value = 7.036731597647719
$ws.Range("D63").Value = 0.4206604806503861
$ws.Range("A64").Value = 62
$ws.Range("B64").Value = 84.50333346370068
$ws.Range("C64").Value = 2.3239409830327
$ws.Range("D64").Value = 0.1503785840451847
$ws.Range("A65").Value = 63
$ws.Range("B65").Value = 85.65725773224199
$ws.Range("C65").Value = 1.993166014106106
$ws.Range("D65").Value = 0.184996217630928
$ws.Range("A66").Value = 64
$ws.Range("B66").Value = 89.2709543434259
$ws.Range("C66").Value = 8.235744315598357
$ws.Range("D66").Value = 0.373639493822613
$ws.Range("A67").Value = 65
$ws.Range("B67").Value = 89.45417966993168
$ws.Range("C67").Value = 4.624621999651204
$ws.Range("D67").Value = 0.1499262540199049
$ws.Range("A68").Value = 66
$ws.Range("B68").Value = 91.15186981287195
$ws.Range("C68").Value = 11.12389251957844
$ws.Range("D68").Value = 0.4839930543577219
$ws.Range("A69").Value = 67
$ws.Range("B69").Value = 92.5893845031046
$ws.Range("C69").Value = 4.113366479032526
$ws.Range("D69").Value = 0.3674693278165844
$ws.Range("A70").Value = 68
$ws.Range("B70").Value = 94.71002842921969
$ws.Range("C70").Value = 12.93960394383676
$ws.Range("D70").Value = 0.3001748441411214
$ws.Range("A71").Value = 69
$ws.Range("B71").Value = 96.48920304198303
$ws.Range("C71").Value = 12.68272858677714
$ws.Range("D71").Value = 0.4978761357404759
$ws.Range("A72").Value = 70
$ws.Range("B72").Value = 97.13437737036693
$ws.Range("C72").Value = 10.357868416477
$ws.Range("D72").Value = 0.4899175020057303
$ws.Range("A73").Value = 71
$ws.Range("B73").Value = 97.29753279989633
$ws.Range("C73").Value = 5.559017358906497
$ws.Range("D73").Value = 0.1555420610243804
$ws.Range("A74").Value = 72
$ws.Range("B74").Value = 99.05510070507634
$ws.Range("C74").Value = 7.769626445485784
$ws.Range("D74").Value = 0.1488188506050884
